$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (LCSC part "DEBF33D103ZA3B") changed supplier part. The new LCSC
# catalog number goes into D3 as plain text (no longer hyperlinked), and the
# manufacturer part number now also appears in C3.
$ws.Range("D3").Value = "C263258"
$ws.Range("C3").Value = "CY2103MD1IEF4CC0A8"

# The worksheet only keeps a hyperlink on D3 as a side effect of the old
# part number; since that hyperlink needs to disappear, recreate the
# remaining hyperlinks (D20, D6, D10, E10) after clearing all of them -
# Hyperlinks.Delete() operates on the whole sheet's collection.
$targets = @{
    "D20" = "https://lcsc.com/product-detail/SMD-Crystal-Resonators_ZHUHAI-MONEC-ELEC-ESB16-0000F12M25F_C353715.html"
    "D6"  = "https://lcsc.com/product-detail/Fuses-with-Leads-Through-Hole_Xucheng-Elec-5TE-05002R1BT_C140483.html"
    "D10" = "https://aliexpress.ru/item/32828769503.html?spm=a2g0o.productlist.0.0.1e355ed4nyrFI0&algo_pvid=328fed0c-bd05-4a2b-b1ec-ac97e390c2d6&algo_expid=328fed0c-bd05-4a2b-b1ec-ac97e390c2d6-6&btsid=daa234f5-190c-42e0-9c04-f9e153b4c333&ws_ab_test=searchweb0_0,search"
    "E10" = "https://aliexpress.ru/item/32809198141.html?spm=a2g0o.cart.0.0.4d963c00LcAFEO&mp=1"
}

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D20"), $targets["D20"], "", "", "Cristal")
$ws.Hyperlinks.Add($ws.Range("D6"), $targets["D6"])
$ws.Hyperlinks.Add($ws.Range("D10"), $targets["D10"])
$ws.Hyperlinks.Add($ws.Range("E10"), $targets["E10"])

$ws.Range("C3").Select()
